$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colC = @{
    4  = -2.034793027571991
    5  = 3.880748577052473
    6  = 2.174296999091507
    7  = -1.942476814417471
    8  = 4.926006686093287
    9  = 0.7614971083056643
    10 = 3.004126378695804
    11 = 2.384525276921168
    12 = 2.677874785158663
    13 = 3.746278935679004
    14 = 3.278383090085746
    15 = 2.877357105987888
    16 = -0.4001895765463725
    17 = -2.321114556364801
    18 = -0.9245282159112467
    19 = 0.2384226118222088
}

$colE = @{
    4  = 0.6952816881563351
    5  = 2.936333428994109
    6  = 1.706732094556851
    7  = 1.288975737543607
    8  = 2.643411312704802
    9  = 2.348955682567344
    10 = 2.358460376580607
    11 = 2.349880264276161
    12 = 2.507403033330702
    13 = 1.868164064786093
    14 = 2.401026764575831
    15 = 2.770626214993133
    16 = 2.524513388369543
    17 = 1.637366184014355
    18 = 0.9277865862836965
    19 = 0.8054161303035379
}

foreach ($row in 4..19) {
    $ws.Cells.Item($row, 3).Value = $colC[$row]
    $ws.Cells.Item($row, 5).Value = $colE[$row]
}
